$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Update F column (Max. Slip Frequency [Hz]) for relay rows 2 through 15: 13 -> 2
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = 2
}

# Update C column (Max. Slip Voltage [%]) for the 480V relays (rows 10, 14, 15): 5 -> 10
$ws.Cells.Item(10, 3).Value = 10
$ws.Cells.Item(14, 3).Value = 10
$ws.Cells.Item(15, 3).Value = 10

# Update selection to C12
$ws.Range("C12").Select()
